# Applies the "manual review" audit columns to the fastqFiles sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells G1:K1
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "MANUAL_AUDIT_20200520"
$ws.Range("H1").Value = "FAIL_REASON_20200520"
$ws.Range("I1").Value = "MANUAL_AUDIT_20207230"
$ws.Range("J1").Value = "FAIL_REASON_20200723"
$ws.Range("K1").Value = "NOTES"

# ---------------------------------------------------------------------------
# 2. Manual-audit values for rows 3-35 (row 2 is left untouched)
# ---------------------------------------------------------------------------
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 0

$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 8
$ws.Range("K25").Value = "Fail - borderline expression of unexpected marker was deemed unnacceptable"

$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 4

$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 4
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 256
$ws.Range("K27").Value = "Fail - unexplained pattern of reads that don't match the other ""good"" replicates"
$ws.Range("K27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 124.6

$ws.Range("G28").Value = 1

$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 4

$ws.Range("G30").Value = 0

$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 1

$ws.Range("G32").Value = 0
$ws.Range("G33").Value = 0

$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 3

$ws.Range("G35").Value = 0

# ---------------------------------------------------------------------------
# 3. Header row font: Calibri -> Cambria
# ---------------------------------------------------------------------------
$ws.Range("A1:K1").Font.Name = "Cambria"

# ---------------------------------------------------------------------------
# 4. Column width: 8.83 -> 8.67 (character units)
# ---------------------------------------------------------------------------
$ws.Range("A1:Z1").EntireColumn.ColumnWidth = 7.836666666666667

# ---------------------------------------------------------------------------
# 5. Page margins: left/right 0.7in -> 0.75in, top/bottom 0.75in -> 1in
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72

# ---------------------------------------------------------------------------
# 6. Active selection
# ---------------------------------------------------------------------------
$ws.Range("K27").Select()
